# The sheet originally contained several extra example/"comment" rows
# (rows 3-8) under the header row. Per the commit message ("deleted some
# unnecessary comments"), remove those extra rows, keeping only the
# header (row 1) and the first data row (row 2). Also update the first
# data row's values to "win" / "damage" to match the target state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining data row's values.
$ws.Range("B2").Value = "win"
$ws.Range("C2").Value = "damage"

# Remove the now-unnecessary rows 3 through 8 entirely.
$ws.Range("A3:C8").EntireRow.Delete()
